$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$ws.Range('D2').NumberFormat = '@'
$ws.Range('D2').Value = '68.734.08'
$ws.Range('E2').Value = '  +2.16%  '
$ws.Range('D3').NumberFormat = '@'
$ws.Range('D3').Value = '3.757.05'
$ws.Range('E3').Value = '  +0.40%  '
$ws.Range('D4').NumberFormat = '@'
$ws.Range('D4').Value = '0.999'
$ws.Range('E4').Value = '  -0.14%  '
$ws.Range('D5').NumberFormat = '@'
$ws.Range('D5').Value = '602.65'
$ws.Range('E5').Value = '  +1.11%  '
$ws.Range('D6').NumberFormat = '@'
$ws.Range('D6').Value = '169.45'
$ws.Range('E6').Value = '  +1.43%  '
$ws.Range('D7').NumberFormat = '@'
$ws.Range('D7').Value = '3.755.23'
$ws.Range('E7').Value = '  +0.36%  '
$ws.Range('E8').Value = '  +0.01%  '
$ws.Range('E9').Value = '  +2.07%  '
$ws.Range('E10').Value = '  -0.98%  '
$ws.Range('D11').NumberFormat = '@'
$ws.Range('D11').Value = '6.33'
$ws.Range('E11').Value = '  +2.45%  '
$ws.Range('D12').NumberFormat = '@'
$ws.Range('D12').Value = '0.464'
$ws.Range('E12').Value = '  +0.39%  '
$ws.Range('D13').NumberFormat = '@'
$ws.Range('D13').Value = '38.37'
$ws.Range('E13').Value = '  +0.62%  '
$ws.Range('E14').Value = '  +0.47%  '
$ws.Range('D15').NumberFormat = '@'
$ws.Range('D15').Value = '4.383.56'
$ws.Range('E15').Value = '  +0.14%  '
$ws.Range('D16').NumberFormat = '@'
$ws.Range('D16').Value = '3.756.18'
$ws.Range('E16').Value = '  +0.08%  '
$ws.Range('D17').NumberFormat = '@'
$ws.Range('D17').Value = '68.738.19'
$ws.Range('E17').Value = '  +2.14%  '
$ws.Range('E18').Value = '  +2.53%  '
$ws.Range('E19').Value = '  +0.18%  '
$ws.Range('E20').Value = '  -0.27%  '
$ws.Range('B21').Value = 'Uniswap'
$ws.Range('C21').Value = 'https://coinranking.com/coin/_H5FVG9iW+uniswap-uni'
$ws.Range('D21').NumberFormat = '@'
$ws.Range('D21').Value = '10.81'
$ws.Range('E21').Value = '  +17.46%  '
$ws.Range('B22').Value = 'BitcoinCash'
$ws.Range('C22').Value = 'https://coinranking.com/coin/ZlZpzOJo43mIo+bitcoincash-bch'
$ws.Range('D22').NumberFormat = '@'
$ws.Range('D22').Value = '496.88'
$ws.Range('E22').Value = '  +1.34%  '
$ws.Range('E23').Value = '  +0.27%  '
$ws.Range('D24').NumberFormat = '@'
$ws.Range('D24').Value = '85.67'
$ws.Range('E24').Value = '  +0.73%  '
$ws.Range('D25').NumberFormat = '@'
$ws.Range('D25').Value = '0.0000150'
$ws.Range('E25').Value = '  +4.03%  '
$ws.Range('E26').Value = '  -0.57%  '
$ws.Range('D27').NumberFormat = '@'
$ws.Range('D27').Value = '12.37'
$ws.Range('E27').Value = '  +1.73%  '
$ws.Range('D28').NumberFormat = '@'
$ws.Range('D28').Value = '10.25'
$ws.Range('E28').Value = '  +2.28%  '
$ws.Range('E29').Value = '  +0.22%  '
$ws.Range('E30').Value = '  +5.53%  '
$ws.Range('E31').Value = '  +1.34%  '
$ws.Range('E32').Value = '  +2.89%  '
$ws.Range('D33').NumberFormat = '@'
$ws.Range('D33').Value = '32.14'
$ws.Range('E33').Value = '  +0.61%  '
$ws.Range('D34').NumberFormat = '@'
$ws.Range('D34').Value = '3.903.66'
$ws.Range('E34').Value = '  +0.34%  '
$ws.Range('D35').NumberFormat = '@'
$ws.Range('D35').Value = '3.689.76'
$ws.Range('E35').Value = '  +0.21%  '
$ws.Range('D36').NumberFormat = '@'
$ws.Range('D36').Value = '0.109'
$ws.Range('E36').Value = '  +0.83%  '
$ws.Range('D37').NumberFormat = '@'
$ws.Range('D37').Value = '0.999'
$ws.Range('E37').Value = '  -0.06%  '
$ws.Range('E38').Value = '  +1.36%  '
$ws.Range('D39').NumberFormat = '@'
$ws.Range('D39').Value = '5.85'
$ws.Range('E39').Value = '  +0.93%  '
$ws.Range('E40').Value = '  -0.28%  '
$ws.Range('D41').NumberFormat = '@'
$ws.Range('D41').Value = '0.326'
$ws.Range('E41').Value = '  +0.31%  '
$ws.Range('D42').NumberFormat = '@'
$ws.Range('D42').Value = '446.37'
$ws.Range('E42').Value = '  -3.02%  '
$ws.Range('D43').NumberFormat = '@'
$ws.Range('D43').Value = '48.88'
$ws.Range('E43').Value = '  -0.45%  '
$ws.Range('E44').Value = '  +0.37%  '
$ws.Range('D45').NumberFormat = '@'
$ws.Range('D45').Value = '2.88'
$ws.Range('E45').Value = '  +1.29%  '
$ws.Range('D46').NumberFormat = '@'
$ws.Range('D46').Value = '8.51'
$ws.Range('E46').Value = '  +1.94%  '
$ws.Range('E47').Value = '  +0.00%  '
$ws.Range('D48').NumberFormat = '@'
$ws.Range('D48').Value = '40.63'
$ws.Range('E48').Value = '  +1.18%  '
$ws.Range('D49').NumberFormat = '@'
$ws.Range('D49').Value = '2.850.14'
$ws.Range('E49').Value = '  +1.67%  '
$ws.Range('D50').NumberFormat = '@'
$ws.Range('D50').Value = '142.28'
$ws.Range('E50').Value = '  +1.21%  '
$ws.Range('E51').Value = '  +2.67%  '
